# Normalize the "Recorded By" (column G) values on the Session Analysis
# Results sheet: for several multi-author cells, "System"/"system" is
# moved from the front of the comma-separated list to the back (and a
# couple of two-address cells have their order swapped). Single-value
# cells, and cells already ending in "System", are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, admin@admin.com"              = "admin@admin.com, System"
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, System, system"  = "backup@backdoor.com, system, System"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
